$d = $word.ActiveDocument

$replacements = @(
    @{old="889÷2="; new="571÷5="},
    @{old="901÷3="; new="856÷9="},
    @{old="360÷3="; new="218÷4="},
    @{old="229÷4="; new="922÷3="},
    @{old="471÷9="; new="935÷5="},
    @{old="590÷6="; new="609÷2="},
    @{old="148÷8="; new="174÷5="},
    @{old="821÷3="; new="784÷2="},
    @{old="608÷9="; new="176÷7="},
    @{old="876÷5="; new="175÷3="},
    @{old="101÷8="; new="151÷7="},
    @{old="198÷7="; new="548÷7="},
    @{old="108÷8="; new="663÷9="},
    @{old="541÷5="; new="143÷4="},
    @{old="213÷7="; new="971÷3="},
    @{old="823÷6="; new="177÷6="},
    @{old="175÷9="; new="662÷3="},
    @{old="271÷5="; new="103÷7="},
    @{old="860÷4="; new="252÷3="},
    @{old="452÷6="; new="188÷6="},
    @{old="493÷6="; new="953÷5="},
    @{old="698÷3="; new="207÷7="},
    @{old="173÷8="; new="844÷9="},
    @{old="321÷7="; new="601÷8="},
    @{old="600÷8="; new="680÷4="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}
